$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 163910
$ws.Range("C4").Value = 154895
$ws.Range("C5").Value = 9015
$ws.Range("C8").Value = 64.95
